$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header in C1: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# Update C2: trainingaudio/06_titoka3.wav -> train1P2
$ws.Range("C2").Value = "train1P2"

# Update C3: trainingaudio/02_pitito3.wav -> train1P2
$ws.Range("C3").Value = "train1P2"
